$d = $word.ActiveDocument

$para = $d.Paragraphs.Last
$rng = $para.Range
$rng.Text = "Hi."
$rng.LanguageIDFarEast = "zh-CN"
